$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 787.2
$ws.Range("I32").Value = 750
$ws.Range("J32").Value = 812
$ws.Range("K32").Value = 750
$ws.Range("L32").Value = 812
$ws.Range("M32").Value = -424
$ws.Range("N32").Value = -1464

$ws.Range("H40").Value = 2774.1667
$ws.Range("I40").Value = 4280
$ws.Range("J40").Value = 1698.5714
$ws.Range("K40").Value = 4280
$ws.Range("L40").Value = 1698.5714
$ws.Range("M40").Value = -4105
$ws.Range("N40").Value = -2048.5714

$ws.Range("H62").Value = 69678.60000000001
$ws.Range("I62").Value = 202337.6
$ws.Range("J62").Value = 3349.1
$ws.Range("K62").Value = 202337.6
$ws.Range("L62").Value = 3349.1
$ws.Range("M62").Value = -201713.6
$ws.Range("N62").Value = -4597.1

$ws.Range("H64").Value = 3977
$ws.Range("I64").Value = 3500.1177
$ws.Range("K64").Value = 3500.1177
$ws.Range("M64").Value = -3252.1177

$ws.Range("H65").Value = 69678.60000000001
$ws.Range("I65").Value = 202337.6
$ws.Range("J65").Value = 3349.1
$ws.Range("K65").Value = 1011688
$ws.Range("L65").Value = 16745.5
$ws.Range("M65").Value = -1008568
$ws.Range("N65").Value = -22985.5

$ws.Range("H67").Value = 3977
$ws.Range("I67").Value = 3500.1177
$ws.Range("K67").Value = 3500.1177
$ws.Range("M67").Value = -2642.1177

$ws.Range("H132").Value = 1561.88
$ws.Range("I132").Value = 1402.0889
$ws.Range("K132").Value = 4206.2667
$ws.Range("M132").Value = -1676.2667

$ws.Range("H137").Value = 3874.7144
$ws.Range("I137").Value = 2141.4707
$ws.Range("J137").Value = 5511.6665
$ws.Range("K137").Value = 6424.4121
$ws.Range("L137").Value = 16534.9995
$ws.Range("M137").Value = -3874.4121
$ws.Range("N137").Value = -21634.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1505.1666
$ws.Range("I2").Value = 1559.762
$ws.Range("J2").Value = 1377.7778
$ws.Range("K2").Value = 1559.762
$ws.Range("L2").Value = 1377.7778
$ws.Range("M2").Value = -1446.762
$ws.Range("N2").Value = -1603.7778

$ws.Range("H32").Value = 7384.2407
$ws.Range("I32").Value = 5730.8906
$ws.Range("K32").Value = 5730.8906
$ws.Range("M32").Value = -5443.8906

$ws.Range("H61").Value = 6289.1387
$ws.Range("I61").Value = 4011.76
$ws.Range("J61").Value = 11465
$ws.Range("K61").Value = 4011.76
$ws.Range("L61").Value = 11465
$ws.Range("M61").Value = -3799.76
$ws.Range("N61").Value = -11889

$ws.Range("H74").Value = 97688.60000000001
$ws.Range("I74").Value = 112508.31
$ws.Range("J74").Value = 18650.166
$ws.Range("K74").Value = 112508.31
$ws.Range("L74").Value = 18650.166
$ws.Range("M74").Value = -111634.31
$ws.Range("N74").Value = -20398.166

$ws.Range("H77").Value = 97688.60000000001
$ws.Range("I77").Value = 112508.31
$ws.Range("J77").Value = 18650.166
$ws.Range("K77").Value = 562541.55
$ws.Range("L77").Value = 93250.83
$ws.Range("M77").Value = -558173.55
$ws.Range("N77").Value = -101986.83

$ws.Range("H116").Value = 1505.1666
$ws.Range("I116").Value = 1559.762
$ws.Range("J116").Value = 1377.7778
$ws.Range("K116").Value = 1559.762
$ws.Range("L116").Value = 1377.7778
$ws.Range("M116").Value = 734.2380000000001
$ws.Range("N116").Value = -5965.7778

$ws.Range("H132").Value = 3885.0352
$ws.Range("I132").Value = 1146.6765
$ws.Range("J132").Value = 7933.0435
$ws.Range("K132").Value = 3440.0295
$ws.Range("L132").Value = 23799.1305
$ws.Range("M132").Value = -910.0295000000001
$ws.Range("N132").Value = -28859.1305

$ws.Range("H136").Value = 6289.1387
$ws.Range("I136").Value = 4011.76
$ws.Range("J136").Value = 11465
$ws.Range("K136").Value = 12035.28
$ws.Range("L136").Value = 34395
$ws.Range("M136").Value = -9485.280000000001
$ws.Range("N136").Value = -39495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1505.1666
$ws.Range("I3").Value = 1559.762
$ws.Range("J3").Value = 1377.7778
$ws.Range("K3").Value = 1559.762
$ws.Range("L3").Value = 1377.7778
$ws.Range("M3").Value = -1445.762
$ws.Range("N3").Value = -1605.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3091.158
$ws.Range("I31").Value = 2267.3076
$ws.Range("J31").Value = 4876.1665
$ws.Range("K31").Value = 2267.3076
$ws.Range("L31").Value = 4876.1665
$ws.Range("M31").Value = -1972.3076
$ws.Range("N31").Value = -5466.1665

$ws.Range("H34").Value = 3091.158
$ws.Range("I34").Value = 2267.3076
$ws.Range("J34").Value = 4876.1665
$ws.Range("K34").Value = 2267.3076
$ws.Range("L34").Value = 4876.1665
$ws.Range("M34").Value = -2065.3076
$ws.Range("N34").Value = -5280.1665

$ws.Range("H99").Value = 1780.7307
$ws.Range("I99").Value = 1907.6154
$ws.Range("J99").Value = 1653.8462
$ws.Range("K99").Value = 1907.6154
$ws.Range("L99").Value = 1653.8462
$ws.Range("M99").Value = -409.6153999999999
$ws.Range("N99").Value = -4649.8462

$ws.Range("H126").Value = 1780.7307
$ws.Range("I126").Value = 1907.6154
$ws.Range("J126").Value = 1653.8462
$ws.Range("K126").Value = 5722.8462
$ws.Range("L126").Value = 4961.5386
$ws.Range("M126").Value = -3252.8462
$ws.Range("N126").Value = -9901.5386

$ws.Range("H132").Value = 2593.7036
$ws.Range("I132").Value = 1967.6364
$ws.Range("J132").Value = 5348.4
$ws.Range("K132").Value = 5902.9092
$ws.Range("L132").Value = 16045.2
$ws.Range("M132").Value = -3372.9092
$ws.Range("N132").Value = -21105.2

$ws.Range("H134").Value = 31275.756
$ws.Range("I134").Value = 114744.78
$ws.Range("J134").Value = 4446.4287
$ws.Range("K134").Value = 344234.34
$ws.Range("L134").Value = 13339.2861
$ws.Range("M134").Value = -341699.34
$ws.Range("N134").Value = -18409.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 10000
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H103").Value = 29000
$ws.Range("J103").Value = 29000
$ws.Range("L103").Value = 29000
$ws.Range("N103").Value = -31344

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H132").Value = 5431.0356
$ws.Range("I132").Value = 11447.5
$ws.Range("K132").Value = 34342.5
$ws.Range("M132").Value = -31812.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4096.875
$ws.Range("I7").Value = 4000
$ws.Range("K7").Value = 4000
$ws.Range("M7").Value = -3888

$ws.Range("H23").Value = 258125.25
$ws.Range("I23").Value = 337500
$ws.Range("J23").Value = 20001
$ws.Range("K23").Value = 337500
$ws.Range("L23").Value = 20001
$ws.Range("M23").Value = -337270
$ws.Range("N23").Value = -20461

$ws.Range("H24").Value = 14600
$ws.Range("J24").Value = 14600
$ws.Range("L24").Value = 14600
$ws.Range("N24").Value = -15286

$ws.Range("H25").Value = 19698.666
$ws.Range("I25").Value = 6766.6665
$ws.Range("J25").Value = 26164.666
$ws.Range("K25").Value = 6766.6665
$ws.Range("L25").Value = 26164.666
$ws.Range("M25").Value = -6536.6665
$ws.Range("N25").Value = -26624.666

$ws.Range("H46").Value = 790.6923
$ws.Range("I46").Value = 607.5
$ws.Range("J46").Value = 872.1111
$ws.Range("K46").Value = 607.5
$ws.Range("L46").Value = 872.1111
$ws.Range("M46").Value = -419.5
$ws.Range("N46").Value = -1248.1111

$ws.Range("H55").Value = 200691.05
$ws.Range("I55").Value = 364358.28
$ws.Range("K55").Value = 364358.28
$ws.Range("M55").Value = -364185.28

$ws.Range("H126").Value = 4096.875
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 10000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -20826

$ws.Range("H42").Value = 52536.75
$ws.Range("J42").Value = 52536.75
$ws.Range("L42").Value = 52536.75
$ws.Range("N42").Value = -53292.75

$ws.Range("H94").Value = 98000
$ws.Range("J94").Value = 98000
$ws.Range("L94").Value = 98000
$ws.Range("N94").Value = -99802

$ws.Range("H132").Value = 2505.568
$ws.Range("I132").Value = 2232.375
$ws.Range("J132").Value = 2833.4
$ws.Range("K132").Value = 6697.125
$ws.Range("L132").Value = 8500.200000000001
$ws.Range("M132").Value = -4167.125
$ws.Range("N132").Value = -13560.2
